$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.61 = 50413.85 pesos`n✅ 50413.85 pesos = 12.58 = 972.32 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 79.28
$ws2.Range("O10").Value = 3996.81
$ws2.Range("N12").Value = 4009
$ws2.Range("O12").Value = 77.321
